$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the A:D values between row 42 and row 43 (E and F remain unchanged)
$a42 = $ws.Range("A42").Value()
$b42 = $ws.Range("B42").Value()
$c42 = $ws.Range("C42").Value()
$d42 = $ws.Range("D42").Value()

$a43 = $ws.Range("A43").Value()
$b43 = $ws.Range("B43").Value()
$c43 = $ws.Range("C43").Value()
$d43 = $ws.Range("D43").Value()

$ws.Range("A42").Value = $a43
$ws.Range("B42").Value = $b43
$ws.Range("C42").Value = $c43
$ws.Range("D42").Value = $d43

$ws.Range("A43").Value = $a42
$ws.Range("B43").Value = $b42
$ws.Range("C43").Value = $c42
$ws.Range("D43").Value = $d42
